# Refactor code to save results in a specified folder.
# Re-running the sliding-window model produced new "IPC PO" predictions
# (column C) for rows 2..51; the dependent DELTA (column D), DELTA^2
# (column E) and the TOTAL / MSE summary cells are updated to match the
# newly produced model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number => @(IPC PO, DELTA, DELTA^2) for the new run, rows 2..51
$data = @{
    2 = @(27.82032490848616, -1.329675091513835, 1.768035848992325)
    3 = @(29.28318637371153, -0.06681362628847154, 0.004464060657815535)
    4 = @(29.00570591093511, -0.3642940890648916, 0.1327101833276192)
    5 = @(28.71038135676348, -0.8296186432365218, 0.6882670932056072)
    6 = @(28.44399932177484, -1.106000678225158, 1.223237500234509)
    7 = @(29.22256232725998, -0.5274376727400245, 0.2781904986254132)
    8 = @(30.09548599497818, 0.2554859949781836, 0.06527309362999245)
    9 = @(29.47642882762488, -0.3335711723751231, 0.1112697270397141)
    10 = @(30.680277994943, 0.7602779949429959, 0.5780226295945421)
    11 = @(30.08983849801517, 0.109838498015165, 0.0120644956462274)
    12 = @(30.0804113654337, 0.04041136543370527, 0.001633078456216469)
    13 = @(29.95985272288711, -0.2501472771128945, 0.06257366024699525)
    14 = @(29.73545450686738, -0.4845454931326145, 0.2347843349151286)
    15 = @(30.37821449262295, -0.001785507377050521, 0.000003188036593501833)
    16 = @(30.54320638908418, 0.1032063890841819, 0.01065155874779554)
    17 = @(30.56629890651797, 0.08629890651796757, 0.007447501266196905)
    18 = @(31.37525106180708, 0.6852510618070795, 0.4695690177077299)
    19 = @(30.89061187403108, 0.1406118740310838, 0.01977169911853337)
    20 = @(31.83861355681801, 0.8986135568180096, 0.8075063244971141)
    21 = @(30.88580861726953, -0.06419138273046698, 0.004120533616849295)
    22 = @(31.7919526936309, 0.7719526936308974, 0.5959109612039981)
    23 = @(31.99901322385107, 0.8790132238510644, 0.7726642477050415)
    24 = @(32.58663422909679, 1.306634229096794, 1.707293008647373)
    25 = @(31.99138543685247, 0.6113854368524692, 0.3737921523952846)
    26 = @(32.43089861652436, 0.8508986165243613, 0.724028455603072)
    27 = @(31.90531148879703, 0.2553114887970267, 0.06518395631175429)
    28 = @(32.90829830425423, 1.028298304254232, 1.057397402532129)
    29 = @(32.90915359681569, 0.6291535968156907, 0.3958342483861207)
    30 = @(32.40525259128498, -0.04474740871502547, 0.002002330586709537)
    31 = @(33.92527120739143, 1.075271207391424, 1.156208169445011)
    32 = @(32.69907304263058, -0.2009269573694183, 0.04037164219773206)
    33 = @(32.91411211806948, -0.1858878819305261, 0.0345543046486172)
    34 = @(33.5396928583843, 0.1396928583843007, 0.01951409468357629)
    35 = @(33.85946444839671, 0.1594644483967045, 0.02542891030246523)
    36 = @(34.55092894464192, 0.4509289446419231, 0.2033369131158786)
    37 = @(34.44573863409703, 0.04573863409703449, 0.002092022649062406)
    38 = @(35.12841259634422, 0.2284125963442207, 0.05217231416870792)
    39 = @(35.11275519830228, -0.187244801697716, 0.03506061576281697)
    40 = @(35.507262394196, -0.1927376058040053, 0.03714778469106016)
    41 = @(35.75154825222906, -0.5484517477709403, 0.3007993196329991)
    42 = @(36.15734149437066, -0.6426585056293419, 0.4130099548577389)
    43 = @(36.83780455716139, -0.4621954428386061, 0.2136246273807752)
    44 = @(38.4470830984661, 0.5470830984661035, 0.2992999166272723)
    45 = @(38.71444306322609, 0.2144430632260921, 0.04598582736578972)
    46 = @(38.92065921215465, 0.02065921215465494, 0.0004268030468510423)
    47 = @(39.69124988482131, 0.2912498848213119, 0.08482649540842743)
    48 = @(40.45294454991505, 0.552944549915054, 0.3057476752807617)
    49 = @(38.67116111609547, -1.428838883904533, 2.04158055615755)
    50 = @(40.01416341001176, -0.5858365899882401, 0.3432045101690493)
    51 = @(39.64829065710045, -1.251709342899545, 1.566776279102011)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value2 = $vals[0]   # IPC PO
    $ws.Cells.Item($row, 4).Value2 = $vals[1]   # DELTA
    $ws.Cells.Item($row, 5).Value2 = $vals[2]   # DELTA^2
}

# TOTAL row (52): sum of DELTA -> column C, sum of DELTA^2 -> column E
$ws.Cells.Item(52, 3).Value2 = 2.049215926944782
$ws.Cells.Item(52, 5).Value2 = 19.39487152762855

# MSE row (53): mean of DELTA^2 -> column E
$ws.Cells.Item(53, 5).Value2 = 0.387897430552571
